$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto price/volume data as captured in the latest GitHub Actions run.
# Force text number format on every touched cell so values such as "109.40" or
# "51.919.44" are preserved exactly as text (matching the source inlineStr cells)
# instead of being auto-converted to numbers by Excel.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "51.919.44"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.41%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.788.70"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "360.91"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.05%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "109.40"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -3.60%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.559"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -3.05%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -1.79%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.10"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -3.49%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -1.64%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +1.05%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "19.51"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -2.46%  "
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -2.82%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.227.53"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -2.05%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.792.42"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -2.27%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.938"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "51.893.88"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.21%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.46"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -1.63%  "
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -2.34%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.10"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -3.58%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0₃0975"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -1.79%  "
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.05%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "269.62"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.65%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -2.41%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.52"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -2.49%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.07%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.160"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +14.13%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "10.28"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -1.28%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.26"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.31%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0469"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +1.98%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "51.99"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -3.29%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "34.24"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +0.17%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -2.74%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.11%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.23"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -3.41%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.01%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "19.03"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +3.75%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -2.15%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -3.98%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.61"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +2.15%  "
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -2.11%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.25"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -1.19%  "
$ws.Range("B44").NumberFormat = "@"
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").NumberFormat = "@"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "21.99"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -7.49%  "
$ws.Range("B45").NumberFormat = "@"
$ws.Range("B45").Value = "Monero"
$ws.Range("C45").NumberFormat = "@"
$ws.Range("C45").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "119.13"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -7.29%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.083.05"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -1.39%  "
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -4.13%  "
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -1.87%  "
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -1.12%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.951"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -5.05%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.85"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -2.45%  "

Write-Host "Applied 83 cell updates"
